# Fix macrons lost in transfer from Word to CSV for the "Maori Name" codeset entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Māori Name"
$ws.Range("C4").Value = "A name for the organisation in te reo Māori"
